$d = $word.ActiveDocument

function Merge-ParaText($para) {
    # Replaces a paragraph's text with itself via Find/Replace, which causes
    # the COM host to collapse any split runs / proofErr tags into a single run.
    # Paragraph.Range.Text carries a trailing paragraph-mark character (\r) that
    # Find.Execute won't match literally, so trim it off first.
    $t = $para.Range.Text
    $t = $t.Substring(0, $t.Length - 1)
    $para.Range.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}

# --- 1) Strip spell-check run-splitting (proofErr) by re-writing paragraph text in place ---

# "The ToDo list must be connected through a front end HTML interface, using a PHP and MySql backend for data storage."
Merge-ParaText $d.Paragraphs(6)

# "The ToDo list should be able to add, view, and delete tasks."
Merge-ParaText $d.Paragraphs(10)

# "The deployment of this list is through the Apache server using PHP and MySql"
Merge-ParaText $d.Paragraphs(11)

# "storeData class, that would have the code pushed to the database and stored."
Merge-ParaText $d.Paragraphs(42)

# "pushData class, which would push the data from the PHP query to the SQL database"
Merge-ParaText $d.Paragraphs(43)

# "arrangeTasks class, which would be more so for the user front end, allowing the user to move tasks up or down in the list and rank them by priority."
Merge-ParaText $d.Paragraphs(44)

# --- 2) Rework the "Test Case" section ---
# Before:  Test Case / <empty> / <empty> / Name of task: Learning how to use XAMPP...
# After:   Test Case / <new summary paragraph w/ _GoBack bookmark> / <empty> / Name of task: ...

$firstBlank = $d.Paragraphs(46)
$secondBlank = $d.Paragraphs(47)

# Put the bookmark on the (still empty) paragraph BEFORE inserting text, so that
# once text is typed in, the bookmark settles right after the inserted run
# (matching Word's normal behaviour for the "_GoBack" last-edit bookmark).
$bmRange = $firstBlank.Range.Duplicate
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$firstBlank.Range.Text = "I was unable to use the add/delete buttons on my front end interface, but was able to add multiple hard coded test cases using the PHP code in my connect.php file."

# Remove the now-redundant second blank paragraph so a single blank line remains
# between the new paragraph and "Name of task: Learning how to use XAMPP...".
$secondBlank.Range.Delete() | Out-Null

# Re-merge the runs in the following "Name of task" / "Notes" paragraphs that had
# proofErr-driven run splits (indices shifted by -1 vs. the original doc because
# we removed one blank paragraph above).

# "Name of task: Learning how to use XAMPP and myPHPAdmin"
Merge-ParaText $d.Paragraphs(48)

# "Notes: Had to learn how to use both programs, and add the data to the databases through the php and MySQL coding"
Merge-ParaText $d.Paragraphs(51)

# "Name of task: Learning PHP and mySql" (also drops the stray empty leading run)
Merge-ParaText $d.Paragraphs(54)

# "Notes: Using an HTML interface for the user and a PHP and mySql backend, I was able to connect to the database, allowing me to store and view information that was used in my code."
Merge-ParaText $d.Paragraphs(57)

# "Notes:  Created detailed information and diagrams about how my todo list application was setup."
Merge-ParaText $d.Paragraphs(63)
